# Auto-generated edit script: apply numeric corrections to Lamia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1177.3522
$ws.Range("I15").Value = 1177.3522
$ws.Range("K15").Value = 3532.0566
$ws.Range("M15").Value = -3363.0566
$ws.Range("H18").Value = 126.5
$ws.Range("I18").Value = 128
$ws.Range("K18").Value = 128
$ws.Range("M18").Value = 156
$ws.Range("H80").Value = 3284.6875
$ws.Range("J80").Value = 2640
$ws.Range("L80").Value = 7920
$ws.Range("N80").Value = -9916
$ws.Range("H83").Value = 3284.6875
$ws.Range("J83").Value = 2640
$ws.Range("L83").Value = 23760
$ws.Range("N83").Value = -33744
$ws.Range("H92").Value = 5276.273
$ws.Range("I92").Value = 3991.8333
$ws.Range("K92").Value = 3991.8333
$ws.Range("M92").Value = -2743.8333
$ws.Range("H95").Value = 62852
$ws.Range("J95").Value = 62852
$ws.Range("L95").Value = 62852
$ws.Range("N95").Value = -68344
$ws.Range("H106").Value = 3130.5
$ws.Range("I106").Value = 1809.909
$ws.Range("K106").Value = 1809.909
$ws.Range("M106").Value = -1178.909
$ws.Range("H113").Value = 11624
$ws.Range("I113").Value = 12965.143
$ws.Range("K113").Value = 12965.143
$ws.Range("M113").Value = -9711.143
$ws.Range("H125").Value = 1806.4
$ws.Range("I125").Value = 1365.125
$ws.Range("K125").Value = 12286.125
$ws.Range("M125").Value = -9826.125
$ws.Range("H138").Value = 3334.802
$ws.Range("I138").Value = 1336.6666
$ws.Range("J138").Value = 3795.9102
$ws.Range("K138").Value = 4009.9998
$ws.Range("L138").Value = 11387.7306
$ws.Range("M138").Value = 1130.0002
$ws.Range("N138").Value = -21667.7306

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4429.3
$ws.Range("J2").Value = 9994.454
$ws.Range("L2").Value = 9994.454
$ws.Range("N2").Value = -10220.454
$ws.Range("H8").Value = 5000000
$ws.Range("I8").Value = 5000000
$ws.Range("K8").Value = 5000000
$ws.Range("M8").Value = -4999856
$ws.Range("H35").Value = 12500
$ws.Range("I35").Value = 12500
$ws.Range("K35").Value = 12500
$ws.Range("M35").Value = -12094
$ws.Range("H97").Value = 933.5135
$ws.Range("I97").Value = 809.069
$ws.Range("J97").Value = 1384.625
$ws.Range("K97").Value = 809.069
$ws.Range("L97").Value = 1384.625
$ws.Range("M97").Value = -313.069
$ws.Range("N97").Value = -2376.625
$ws.Range("H116").Value = 4429.3
$ws.Range("J116").Value = 9994.454
$ws.Range("L116").Value = 9994.454
$ws.Range("N116").Value = -14582.454
$ws.Range("H121").Value = 64234
$ws.Range("J121").Value = 64234
$ws.Range("L121").Value = 64234
$ws.Range("N121").Value = -67728
$ws.Range("H133").Value = 99980.5
$ws.Range("J133").Value = 99980.5
$ws.Range("L133").Value = 99980.5
$ws.Range("N133").Value = -105040.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4429.3
$ws.Range("J3").Value = 9994.454
$ws.Range("L3").Value = 9994.454
$ws.Range("N3").Value = -10222.454
$ws.Range("H31").Value = 9923
$ws.Range("I31").Value = 9923
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 9923
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -9671
$ws.Range("N31").ClearContents()
$ws.Range("H86").Value = 3754.8076
$ws.Range("I86").Value = 3127.3157
$ws.Range("J86").Value = 5458
$ws.Range("K86").Value = 3127.3157
$ws.Range("L86").Value = 5458
$ws.Range("M86").Value = -2004.3157
$ws.Range("N86").Value = -7704
$ws.Range("H89").Value = 3754.8076
$ws.Range("I89").Value = 3127.3157
$ws.Range("J89").Value = 5458
$ws.Range("K89").Value = 15636.5785
$ws.Range("L89").Value = 27290
$ws.Range("M89").Value = -10020.5785
$ws.Range("N89").Value = -38522
$ws.Range("H94").Value = 1239.65
$ws.Range("I94").Value = 1367.7646
$ws.Range("J94").Value = 513.6667
$ws.Range("K94").Value = 1367.7646
$ws.Range("L94").Value = 513.6667
$ws.Range("M94").Value = -916.7646
$ws.Range("N94").Value = -1415.6667
$ws.Range("H99").Value = 1929.6
$ws.Range("I99").Value = 1772.4445
$ws.Range("J99").Value = 2165.3333
$ws.Range("K99").Value = 1772.4445
$ws.Range("L99").Value = 2165.3333
$ws.Range("M99").Value = -274.4445000000001
$ws.Range("N99").Value = -5161.3333
$ws.Range("H134").Value = 3776.9546
$ws.Range("I134").Value = 2971.111
$ws.Range("K134").Value = 8913.332999999999
$ws.Range("M134").Value = -6378.332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5443.778
$ws.Range("I22").Value = 866.6667
$ws.Range("J22").Value = 7732.3335
$ws.Range("K22").Value = 866.6667
$ws.Range("L22").Value = 7732.3335
$ws.Range("M22").Value = -516.6667
$ws.Range("N22").Value = -8432.333500000001
$ws.Range("H31").Value = 45121
$ws.Range("I31").Value = 2236.4707
$ws.Range("K31").Value = 2236.4707
$ws.Range("M31").Value = -1941.4707
$ws.Range("H34").Value = 45121
$ws.Range("I34").Value = 2236.4707
$ws.Range("K34").Value = 2236.4707
$ws.Range("M34").Value = -2034.4707
$ws.Range("H140").Value = 69269.336
$ws.Range("J140").Value = 69269.336
$ws.Range("L140").Value = 69269.336
$ws.Range("N140").Value = -79629.336
$ws.Range("H141").Value = 201158
$ws.Range("J141").Value = 208694.2
$ws.Range("L141").Value = 208694.2
$ws.Range("N141").Value = -219054.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42995.855
$ws.Range("I2").Value = 59.666668
$ws.Range("J2").Value = 54705.727
$ws.Range("K2").Value = 358.000008
$ws.Range("L2").Value = 328234.362
$ws.Range("M2").Value = -245.000008
$ws.Range("N2").Value = -328460.362
$ws.Range("H23").Value = 899.7619
$ws.Range("I23").Value = 311.2857
$ws.Range("J23").Value = 1194
$ws.Range("K23").Value = 933.8571000000001
$ws.Range("L23").Value = 3582
$ws.Range("M23").Value = -698.8571000000001
$ws.Range("N23").Value = -4052
$ws.Range("H35").Value = 5726.6
$ws.Range("J35").Value = 6658.25
$ws.Range("L35").Value = 19974.75
$ws.Range("N35").Value = -20550.75
$ws.Range("H98").Value = 2209
$ws.Range("J98").Value = 2689.7144
$ws.Range("L98").Value = 8069.1432
$ws.Range("N98").Value = -11065.1432
$ws.Range("H107").Value = 453442.03
$ws.Range("J107").Value = 1250676.1
$ws.Range("L107").Value = 3752028.3
$ws.Range("N107").Value = -3755868.3
$ws.Range("H131").Value = 6265428.5
$ws.Range("I131").Value = 13890332
$ws.Range("J131").Value = 4631520.5
$ws.Range("K131").Value = 41670996
$ws.Range("L131").Value = 13894561.5
$ws.Range("M131").Value = -41665956
$ws.Range("N131").Value = -13904641.5
$ws.Range("H132").Value = 5148.6665
$ws.Range("I132").Value = 4124.4614
$ws.Range("J132").Value = 6359.091
$ws.Range("K132").Value = 37120.1526
$ws.Range("L132").Value = 57231.819
$ws.Range("M132").Value = -34590.1526
$ws.Range("N132").Value = -62291.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 68673.875
$ws.Range("J140").Value = 68673.875
$ws.Range("L140").Value = 68673.875
$ws.Range("N140").Value = -79033.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2369.1667
$ws.Range("I22").Value = 1071
$ws.Range("K22").Value = 1071
$ws.Range("M22").Value = -776
$ws.Range("H27").Value = 2369.1667
$ws.Range("I27").Value = 1071
$ws.Range("K27").Value = 1071
$ws.Range("M27").Value = -964
$ws.Range("H111").Value = 73684
$ws.Range("J111").Value = 73684
$ws.Range("L111").Value = 73684
$ws.Range("N111").Value = -81864
$ws.Range("H122").Value = 154957.86
$ws.Range("I122").Value = 166314.28
$ws.Range("K122").Value = 498942.84
$ws.Range("M122").Value = -496492.84
$ws.Range("H132").Value = 7035.391
$ws.Range("I132").Value = 6863.8184
$ws.Range("J132").Value = 7192.6665
$ws.Range("K132").Value = 20591.4552
$ws.Range("L132").Value = 21577.9995
$ws.Range("M132").Value = -18061.4552
$ws.Range("N132").Value = -26637.9995
$ws.Range("H137").Value = 67500
$ws.Range("J137").Value = 67500
$ws.Range("L137").Value = 67500
$ws.Range("N137").Value = -77700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 39051.5
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 39051.5
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H107").Value = 1784.4
$ws.Range("I107").Value = 1639.5333
$ws.Range("K107").Value = 4918.5999
$ws.Range("M107").Value = -2998.5999
$ws.Range("H132").Value = 1942.2683
$ws.Range("I132").Value = 1758.8
$ws.Range("J132").Value = 2442.6365
$ws.Range("K132").Value = 5276.4
$ws.Range("L132").Value = 7327.9095
$ws.Range("M132").Value = -2746.4
$ws.Range("N132").Value = -12387.9095
$ws.Range("H136").Value = 3465.8823
$ws.Range("I136").Value = 1884.6666
$ws.Range("K136").Value = 5653.9998
$ws.Range("M136").Value = -3103.9998
